$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")

$problemText = @'
给定一个整数类型的数组 nums，请编写一个能够返回数组“中心索引”的方法。 
 我们是这样定义数组中心索引的：数组中心索引的左侧所有元素相加的和等于右侧所有元素相加的和。 
 如果数组不存在中心索引，那么我们应该返回 -1。如果数组有多个中心索引，那么我们应该返回最靠近左边的那一个。 
 示例 1:
输入:
nums = [1, 7, 3, 6, 5, 6]
输出: 3
解释: 
索引3 (nums[3] = 6) 的左侧数之和(1 + 7 + 3 = 11)，与右侧数之和(5 + 6 = 11)相等。
同时, 3 也是第一个符合要求的中心索引。
 示例 2:
输入:
nums = [1, 2, 3]
输出: -1
解释: 
数组中不存在满足此条件的中心索引。 
 说明: 
 nums 的长度范围为 [0, 10000]。
 任何一个 nums[i] 将会是一个范围在 [-1000, 1000]的整数。 
 Related Topics 数组
'@

$stepsText = @'
1 根据数学公式推导：left + mid + right = sum => 2 * left + mid = sum
2 计算数组的元素之和
2 迭代数组，结束条件是找到这个中间值或者数组元素迭代完成
3 定义left=nums[0]，用于mid左边数据的累加。
4 从第i(i=0)个位置开始，2*left + num[i]是否等于sum,left = left + num[i]
'@

$keywordsText = @'
累加
连续数列之和
公示计算
'@

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 724
$ws.Cells.Item(15, 3).Value = $problemText
$ws.Cells.Item(15, 4).Value = $stepsText
$ws.Cells.Item(15, 5).Value = $keywordsText
$ws.Cells.Item(15, 6).Value = "O(N)"
$ws.Cells.Item(15, 7).Value = "O(1)"

$ws.Rows.Item(15).RowHeight = 409.6

$ws.Range("D15").Select()
